$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-ASCII helper: subscript 3 (U+2083) used in D20 ("0.0#0891" style ShibaInu price).
$sub3 = [char]0x2083

# Column D holds numeric-looking text (e.g. "151.57", "40.119.70") that must stay as
# literal text rather than being parsed into Excel numbers. Force Text format on the
# whole D data range first, write the values, then restore the original (unstyled)
# "Normal" style so no style index changes linger on the cells.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = "40.084.23"
$ws.Range("D3").Value = "2.235.86"
$ws.Range("D5").Value = "293.16"
$ws.Range("D6").Value = "87.06"
$ws.Range("D10").Value = "31.20"
$ws.Range("D11").Value = "0.0790"
$ws.Range("D12").Value = "46.92"
$ws.Range("D14").Value = "6.41"
$ws.Range("D15").Value = "2.583.89"
$ws.Range("D16").Value = "14.10"
$ws.Range("D17").Value = "2.238.62"
$ws.Range("D18").Value = "0.733"
$ws.Range("D19").Value = "40.033.06"
$ws.Range("D20").Value = "0.0${sub3}0889"
$ws.Range("D21").Value = "11.27"
$ws.Range("D23").Value = "65.87"
$ws.Range("D24").Value = "236.21"
$ws.Range("D28").Value = "22.93"
$ws.Range("D30").Value = "9.33"
$ws.Range("D31").Value = "33.26"
$ws.Range("D32").Value = "151.58"
$ws.Range("D33").Value = "0.999"
$ws.Range("D34").Value = "4.95"
$ws.Range("D37").Value = "16.27"
$ws.Range("D38").Value = "2.82"
$ws.Range("D42").Value = "3.83"
$ws.Range("D43").Value = "2.060.62"
$ws.Range("D44").Value = "18.17"
$ws.Range("D45").Value = "0.0269"
$ws.Range("D48").Value = "2.59"
$ws.Range("D49").Value = "72.28"
$ws.Range("D50").Value = "2.443.88"
$ws.Range("D51").Value = "89.43"

$dRange.Style = "Normal"

# Column E values (percentage change) are padded with spaces so Excel keeps them as text
# automatically; no special handling required.
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("E5").Value = "  -1.54%  "
$ws.Range("E6").Value = "  +4.36%  "
$ws.Range("E7").Value = "  +0.72%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +1.01%  "
$ws.Range("E10").Value = "  +6.23%  "
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("E12").Value = "  -1.75%  "
$ws.Range("E13").Value = "  +1.42%  "
$ws.Range("E14").Value = "  +1.10%  "
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("E16").Value = "  -0.64%  "
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("E19").Value = "  +1.48%  "
$ws.Range("E20").Value = "  +1.22%  "
$ws.Range("E21").Value = "  +8.85%  "
$ws.Range("E22").Value = "  +1.59%  "
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("E24").Value = "  +3.40%  "
$ws.Range("E25").Value = "  -0.06%  "
$ws.Range("E26").Value = "  +2.60%  "
$ws.Range("E27").Value = "  +2.29%  "
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("E29").Value = "  +2.57%  "
$ws.Range("E30").Value = "  +2.03%  "
$ws.Range("E31").Value = "  +3.35%  "
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("E33").Value = "  +0.13%  "
$ws.Range("E34").Value = "  +1.59%  "
$ws.Range("E35").Value = "  +3.54%  "
$ws.Range("E36").Value = "  +2.20%  "
$ws.Range("E37").Value = "  +6.16%  "
$ws.Range("E38").Value = "  +6.31%  "
$ws.Range("E39").Value = "  +1.71%  "
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("E41").Value = "  +3.69%  "
$ws.Range("E42").Value = "  +4.65%  "
$ws.Range("E43").Value = "  +7.71%  "
$ws.Range("E44").Value = "  +11.83%  "
$ws.Range("E45").Value = "  +3.75%  "
$ws.Range("E46").Value = "  +4.38%  "
$ws.Range("E47").Value = "  +8.89%  "
$ws.Range("E48").Value = "  -1.50%  "
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("E51").Value = "  +2.19%  "
